{"js": "// 1) Bump the version number: \"Versi\u00f3n 0.2\" -> \"Versi\u00f3n 0.3\"\n//    The digit \"2\" lives in its own run right after \"Versi\u00f3n 0.\", so scope\n//    the search to the paragraph that contains \"Versi\u00f3n 0.\" and replace the\n//    lone \"2\" token (avoids touching any other \"2\" elsewhere in the doc).\nconst verMatches = context.document.body.search(\"Versi\u00f3n 0.\", { matchCase: true });\nverMatches.load(\"text,items\");\nawait context.sync();\n\nif (verMatches.items.length > 0) {\n  const verParagraphs = verMatches.items[0].paragraphs;\n  verParagraphs.load(\"text\");\n  await context.sync();\n  const verPara = verParagraphs.items[0];\n\n  const digitMatches = verPara.search(\"2\", { matchCase: true });\n  digitMatches.load(\"text,items\");\n  await context.sync();\n  if (digitMatches.items.length > 0) {\n    digitMatches.items[0].insertText(\"3\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 2) Just above the RFC sentence, insert four blank paragraphs and remove\n//    the phrase \"en may\u00fasculas y \" from \"...tecleen su RFC, en may\u00fasculas y\n//    sin homoclave...\" (users must now type it in lower case).\nconst rfcMatches = context.document.body.search(\"tecleen su RFC, en may\u00fasculas y sin homoclave\", { matchCase: true });\nrfcMatches.load(\"text,items\");\nawait context.sync();\n\nif (rfcMatches.items.length > 0) {\n  const rfcParagraphs = rfcMatches.items[0].paragraphs;\n  rfcParagraphs.load(\"text\");\n  await context.sync();\n  const rfcPara = rfcParagraphs.items[0];\n\n  // Insert 4 new empty (Normal-style) paragraphs right before the RFC one.\n  for (let i = 0; i < 4; i++) {\n    rfcPara.insertParagraph(\"\", Word.InsertLocation.before);\n  }\n  await context.sync();\n\n  // Remove \"en may\u00fasculas y \" from the sentence.\n  const phraseMatches = rfcPara.search(\"en may\u00fasculas y \", { matchCase: true });\n  phraseMatches.load(\"text,items\");\n  await context.sync();\n  if (phraseMatches.items.length > 0) {\n    phraseMatches.items[0].insertText(\"\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Applies the \"Nueva versi\u00f3n, 0.5.21, usuario en minusculas.\" edit:\n#   1) Bumps the version number: \"Versi\u00f3n 0.2\" -> \"Versi\u00f3n 0.3\"\n#   2) Adds four blank lines right above the RFC instructions paragraph\n#   3) Drops \"en may\u00fasculas y \" from that sentence (the RFC no longer has\n#      to be typed in upper case)\n\nfunction Get-ParagraphByText($doc, $needle) {\n    $paras = $doc.Paragraphs\n    $count = $paras.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $candidate = $paras.Item($i)\n        if ($candidate.Range.Text.Contains($needle)) {\n            return $candidate\n        }\n    }\n    return $null\n}\n\n$d = $word.ActiveDocument\n\n# --- 1) Version bump: \"Versi\u00f3n 0.2\" -> \"Versi\u00f3n 0.3\" -----------------------\n$verPara = Get-ParagraphByText $d \"Versi\u00f3n 0.\"\nif ($verPara -ne $null) {\n    $verRange = $verPara.Range\n    $verFind = $verRange.Find\n    $verFind.Text = \"2\"\n    $verFound = $verFind.Execute()\n    if ($verFound) {\n        $verRange.Text = \"3\"\n    }\n}\n\n# --- 2 & 3) RFC paragraph: drop \"en may\u00fasculas y \" and add 4 blank lines --\n$rfcPara = Get-ParagraphByText $d \"tecleen su RFC, en may\u00fasculas y sin homoclave\"\nif ($rfcPara -ne $null) {\n    $rfcRange = $rfcPara.Range\n    $rfcFind = $rfcRange.Find\n    $rfcFind.Text = \"en may\u00fasculas y \"\n    $rfcFind.Replacement.Text = \"\"\n    $rfcFind.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n    # Re-resolve the (now shorter) paragraph and insert 4 empty paragraphs\n    # immediately before it.\n    $rfcParaAgain = Get-ParagraphByText $d \"tecleen su RFC, sin homoclave\"\n    if ($rfcParaAgain -ne $null) {\n        $rfcRangeAgain = $rfcParaAgain.Range\n        for ($i = 0; $i -lt 4; $i++) {\n            $rfcRangeAgain.InsertParagraphBefore()\n        }\n    }\n}\n"}
